$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INAP")

# Insert two new columns before column D (shifts existing D:K data to F:M),
# then copy the number/date formatting from the (shifted) column F into the
# two freshly inserted columns so D:E inherit the same per-row styles that
# column D used to have before the insert.
$ws.Range("D:E").EntireColumn.Insert()
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null
$ws.Range("E7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Period Ending header rows (dates) ----
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373

# ---- Income statement block ----
$ws.Range("D8").Value = 78200
$ws.Range("E8").Value = 83000

$ws.Range("D9").Value = 33700
$ws.Range("E9").Value = 36900

$ws.Range("D10").Value = 44500
$ws.Range("E10").Value = 46200

$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"

$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

$ws.Range("D14").Value = 2300
$ws.Range("E14").Value = 2300

$ws.Range("D15").Value = 23600
$ws.Range("E15").Value = 23400

$ws.Range("D17").Value = 77300
$ws.Range("E17").Value = 80800

$ws.Range("D18").Value = 900
$ws.Range("E18").Value = 2200

$ws.Range("D20").Value = 300
$ws.Range("E20").Value = -200

$ws.Range("D21").Value = 24800
$ws.Range("E21").Value = 25400

$ws.Range("D22").Value = 20300
$ws.Range("E22").Value = 16900

$ws.Range("D23").Value = -19100
$ws.Range("E23").Value = -14900

$ws.Range("D24").Value = 300
$ws.Range("E24").Value = 200

$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

$ws.Range("D26").Value = -19400
$ws.Range("E26").Value = -15100

$ws.Range("D27").Value = -19400
$ws.Range("E27").Value = -15100

$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0

$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0

$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0

$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0

$ws.Range("D32").Value = -300
$ws.Range("E32").Value = 200

$ws.Range("D33").Value = -19400
$ws.Range("E33").Value = -15100

$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0

$ws.Range("D35").Value = -19400
$ws.Range("E35").Value = -15100

# ---- Balance sheet block ----
$ws.Range("D41").Value = 17800
$ws.Range("E41").Value = 11800

$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0

$ws.Range("D43").Value = 28900
$ws.Range("E43").Value = 31000

$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0

$ws.Range("D45").Value = 7400
$ws.Range("E45").Value = 9500

$ws.Range("D46").Value = 54100
$ws.Range("E46").Value = 52400

$ws.Range("D47").Value = 16100
$ws.Range("E47").Value = 12800

$ws.Range("D48").Value = 478100
$ws.Range("E48").Value = 477400

$ws.Range("D49").Value = 189300
$ws.Range("E49").Value = 191400

$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0

$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0

$ws.Range("D52").Value = 7400
$ws.Range("E52").Value = 12100

$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0

$ws.Range("D54").Value = 744900
$ws.Range("E54").Value = 746000

$ws.Range("D57").Value = 23400
$ws.Range("E57").Value = 32200

$ws.Range("D58").Value = 9400
$ws.Range("E58").Value = 28300

$ws.Range("D59").Value = 27200
$ws.Range("E59").Value = 29500

$ws.Range("D60").Value = 60000
$ws.Range("E60").Value = 90000

$ws.Range("D61").Value = 677700
$ws.Range("E61").Value = 667900

$ws.Range("D62").Value = 7200
$ws.Range("E62").Value = 7100

$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0

$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0

$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0

$ws.Range("D66").Value = 747700
$ws.Range("E66").Value = 767800

$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0

$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0

$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0

$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0

$ws.Range("D72").Value = -1363000
$ws.Range("E72").Value = -1343600

$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0

$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0

$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0

$ws.Range("D76").Value = -2700
$ws.Range("E76").Value = -21800

$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0

# ---- Cash flow block ----
$ws.Range("D81").Value = -19400
$ws.Range("E81").Value = -15100

$ws.Range("D83").Value = 23600
$ws.Range("E83").Value = 23400

$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0

$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0

$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0

$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0

$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0

$ws.Range("D89").Value = 5400
$ws.Range("E89").Value = 10300

# Row 91 was restated across every quarter column (not a pure shift)
$ws.Range("D91").Value = -11000
$ws.Range("E91").Value = -11200
$ws.Range("F91").Value = -10000
$ws.Range("G91").Value = -6100
$ws.Range("H91").Value = -12500
$ws.Range("I91").Value = -10900
$ws.Range("J91").Value = -6500
$ws.Range("K91").Value = -6000
$ws.Range("L91").Value = -6200
$ws.Range("M91").Value = -12900

$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0

$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0

$ws.Range("D94").Value = -12300
$ws.Range("E94").Value = -12000

$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0

$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0

$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0

$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0

$ws.Range("D100").Value = 12900
$ws.Range("E100").Value = -1200

$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0

$ws.Range("D102").Value = 6000
$ws.Range("E102").Value = -2900
